# Update scaling mappings to proper format
#
# The "year" sheet stored multi-year selections (e.g. "2006, 2010") as a
# single text cell. This script splits that into one row per year, each
# holding a proper numeric year value, and carries the rest of the row's
# key columns (iso/inv_sector/pre_ext_year/post_ext_year/post_ext_method)
# along with it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("year")

# Duplicate row 2 (the "grc" / "all" row whose E column held "2006, 2010")
# down into a new row 3, copying formatting/styles along with it.
$ws.Rows.Item(2).Copy()
$ws.Rows.Item(3).Insert()

# Replace the combined text value with the individual numeric years.
$ws.Cells.Item(2, 5).Value = 2006
$ws.Cells.Item(3, 5).Value = 2010

# Reflect where the user's selection ended up after the edit.
$ws.Range("E4").Select()
